$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.694.12'
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.531.35'
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '611.49'
$ws.Range("E5").Value = '  +3.75%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.96'
$ws.Range("E6").Value = '  -1.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -1.51%  '
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.213'
$ws.Range("E9").Value = '  +4.96%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.640'
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.33'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000308'
$ws.Range("E12").Value = '  -1.51%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.42'
$ws.Range("E13").Value = '  -1.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.090.92'
$ws.Range("E14").Value = '  -1.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '592.78'
$ws.Range("E15").Value = '  +5.25%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '69.733.74'
$ws.Range("E16").Value = '  -1.34%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.543.49'
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.61'
$ws.Range("E18").Value = '  +1.04%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '18.84'
$ws.Range("E19").Value = '  -4.58%  '
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.988'
$ws.Range("E21").Value = '  -3.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '17.45'
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.70'
$ws.Range("E23").Value = '  +0.15%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '96.81'
$ws.Range("E24").Value = '  +1.15%  '
$ws.Range("E25").Value = '  -1.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.97'
$ws.Range("E26").Value = '  -0.57%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.94'
$ws.Range("E27").Value = '  -6.48%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.55'
$ws.Range("E28").Value = '  +3.93%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '31.98'
$ws.Range("E29").Value = '  -1.06%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.98'
$ws.Range("E30").Value = '  -5.64%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.11'
$ws.Range("E31").Value = '  -3.29%  '
$ws.Range("E32").Value = '  -1.93%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '63.29'
$ws.Range("B34").Value = 'dogwifhat'
$ws.Range("C34").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.58'
$ws.Range("E34").Value = '  +16.61%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.24'
$ws.Range("E35").Value = '  -4.40%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '533.72'
$ws.Range("E36").Value = '  -5.56%  '
$ws.Range("E37").Value = '  -0.03%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.401'
$ws.Range("E38").Value = '  -5.09%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.99'
$ws.Range("E39").Value = '  -3.12%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.540.56'
$ws.Range("E40").Value = '  +5.83%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0₃0777'
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.52'
$ws.Range("E42").Value = '  +4.16%  '
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("E44").Value = '  +1.68%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.93'
$ws.Range("E45").Value = '  -1.87%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.39'
$ws.Range("E46").Value = '  -4.29%  '
$ws.Range("E47").Value = '  +2.52%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.09'
$ws.Range("E48").Value = '  -3.83%  '
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").Value = '  -6.98%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '133.92'
$ws.Range("E51").Value = '  -2.89%  '
